$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.269.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.291"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0720"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.57%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.100.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.833.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.643"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.312.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0792"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.74%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +4.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.445.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0191"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.962"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.994.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0499"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0124"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.57%  "
